$d = $word.ActiveDocument

$replacements = @(
    @("14×20=280", "31×38=1178"),
    @("28×51=1428", "18×34=612"),
    @("15×99=1485", "83×90=7470"),
    @("47×48=2256", "22×14=308"),
    @("91×77=7007", "74×63=4662"),
    @("89×34=3026", "37×35=1295"),
    @("38×75=2850", "73×46=3358"),
    @("92×34=3128", "86×14=1204"),
    @("30×62=1860", "40×76=3040"),
    @("20×36=720", "36×72=2592"),
    @("39×57=2223", "73×88=6424"),
    @("20×64=1280", "84×54=4536"),
    @("56×78=4368", "25×21=525"),
    @("57×16=912", "36×11=396"),
    @("68×85=5780", "42×36=1512"),
    @("38×53=2014", "48×70=3360"),
    @("93×24=2232", "89×96=8544"),
    @("93×89=8277", "84×77=6468"),
    @("25×61=1525", "21×70=1470"),
    @("57×90=5130", "57×88=5016"),
    @("86×45=3870", "60×42=2520"),
    @("48×58=2784", "53×79=4187"),
    @("87×28=2436", "91×64=5824"),
    @("20×46=920", "17×47=799"),
    @("14×12=168", "64×67=4288")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
